$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the sheet; this pushes all the
# existing city/case/death data down by two rows (old row 1 -> row 3, etc).
$ws.Rows("1:2").Insert()

# Row 1: generic column headers (bold, bordered, centered header row)
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

# Row 2: real column titles
$ws.Range("A2").Value = "municipio deresidencia"
$ws.Range("B2").Value = "Casos"
$ws.Range("C2").Value = "Óbitos"

# Two additional data rows appended at the bottom of the table
$ws.Range("A51").Value = "outros estados"
$ws.Range("B51").Value = 16

$ws.Range("A52").Value = "outros paises"
$ws.Range("B52").Value = 27

# Style the new header row (A1:C1): bold font, thin border all around,
# centered horizontally and aligned to the top vertically.
$hdr = $ws.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
